# #5: property boat&car done
# Fixes the "汽車" (car) sheet (sheet3): turns row 1 into proper column
# headers (matching the other sheets' header convention) and extends the
# data row with the standard metadata columns (property_category .. index)
# that every other sheet in this workbook already carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Row 1: header row -------------------------------------------------
# Currently row 1 mistakenly duplicates row 2's data values. Replace them
# with the proper header labels (name, capacity, owner, register_date,
# register_reason, acquire_value), keeping the existing (bold/bordered)
# header style already applied to these cells.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

# Extend the header row with the common metadata header labels used by
# every other sheet (property_category, category, date, legislator_name,
# legislator_id, source_file, index). Copy the existing header formatting
# onto the new cells first so they match the rest of the row.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ---- Row 2: data row ----------------------------------------------------
# A2:G2 already hold the correct car record values, leave them untouched.
# Extend the row with the standard metadata values (property_category,
# category, date, legislator_name, legislator_id, source_file, index),
# matching the pattern used on every other sheet. Copy the existing data
# formatting first so the new cells match the rest of the row.
$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# Use a formula that evaluates to the literal text so Excel does not
# auto-convert the "yyyy-mm-dd"-looking text into a date serial number.
$ws.Range("J2").Formula = '="2013-12-26"'
$ws.Range("K2").Value = "潘孟安"
$ws.Range("L2").Value = 1376
$ws.Range("M2").Value = "tmpf07c1"
$ws.Range("N2").Value = 34

# Flatten the formula in J2 back down to a plain static value/string.
$ws.Range("H2:N2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4163)
$excel.CutCopyMode = 0
